$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 0.6
$ws.Range("E2").Value = 0.738
$ws.Range("F2").Value = 0.737
$ws.Range("G2").Value = 0.418
$ws.Range("H2").Value = 1.482
$ws.Range("I2").Value = 7

# Row 3
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 0.6
$ws.Range("E3").Value = 1.016
$ws.Range("F3").Value = 0.935
$ws.Range("G3").Value = 0.264
$ws.Range("H3").Value = 1.788
$ws.Range("I3").Value = 7

# Row 4
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 0.654
$ws.Range("E4").Value = 1.012
$ws.Range("F4").Value = 1.065
$ws.Range("G4").Value = 0.234
$ws.Range("H4").Value = 1.585
$ws.Range("I4").Value = 7

# Row 5
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 0.628
$ws.Range("E5").Value = 1.123
$ws.Range("F5").Value = 1.067
$ws.Range("G5").Value = 0.257
$ws.Range("H5").Value = 1.909
$ws.Range("I5").Value = 7
